$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert a new column before column V (22nd column) to make room for
# the new "server_calculation" field, shifting existing columns right.
$ws.Columns.Item(22).Insert()

# Set header value/style for new V1 cell
$ws.Cells.Item(1, 22).Value = "server_calculation"
$ws.Cells.Item(1, 22).Font.Bold = $true

# Update the active selection to match the target view state
$ws.Range("V2").Select()
